$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaByText($needle) {
    $result = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            $result = $p
        }
    }
    return $result
}

function Replace-ParaXml($needle, $xml) {
    $p = Get-ParaByText $needle
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $rng.InsertXML($xml)
}

# 1. Delete the "Only use ADL where 'necessary'." bullet entirely.
$p = Get-ParaByText "*Only use ADL where*"
$p.Range.Delete()

# 2. "Optimize all components..." gains the relocated _GoBack bookmark at its start.
Replace-ParaXml "*Optimize all components, both in terms of interface and implementation*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>Optimize all components, both in terms of interface and implementation.</w:t></w:r>' +
  '</w:p>'
)

# 3. "Investigate whether..." loses the bookmark + trailing-space run; the space
#    is folded into the main run's text instead.
Replace-ParaXml "*Investigate whether it?s feasible*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Investigate whether it' + [char]0x2019 + 's feasible and worthwhile to detect when HadesMem is being used ' + [char]0x2018 + 'in-process' + [char]0x2019 + ' and drop to lower level implementations of certain functions such as Read/Write (using SEH instead of RPM/WPM for example), and also using different APIs (such as using local threads rather than remote threads, etc). </w:t></w:r>' +
  '</w:p>'
)

# 4. "New Modules" loses its lastRenderedPageBreak marker.
Replace-ParaXml "*New Modules*" (
  '<w:p ' + $w + '>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>New Modules</w:t></w:r>' +
  '</w:p>'
)

# 5. "Remote memory 'pool'..." gains the lastRenderedPageBreak marker.
Replace-ParaXml "*Remote memory*pool*to avoid allocating*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Remote memory ' + [char]0x2018 + 'pool' + [char]0x2019 + ' to avoid allocating entire pages for only a few bytes of data.</w:t></w:r>' +
  '</w:p>'
)

# 6. "Transactional hooking." loses its lastRenderedPageBreak marker.
Replace-ParaXml "*Transactional hooking*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Transactional hooking.</w:t></w:r>' +
  '</w:p>'
)

# 7. "Improved relative instruction rebuilding..." gains the lastRenderedPageBreak marker.
Replace-ParaXml "*Improved relative instruction rebuilding*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Improved relative instruction rebuilding (including conditionals).</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> x64 has far more IP relative instructions than x86.</w:t></w:r>' +
  '</w:p>'
)

# 8. "Full support for writing back to PE file..." loses its lastRenderedPageBreak marker.
Replace-ParaXml "*Full support for writing back to PE file*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Full support for writing back to PE file, including automatically performing adjustments where required to fit in new data or remove unnecessary space.</w:t></w:r>' +
  '</w:p>'
)

# 9. "Improve export forwarding code..." gains the lastRenderedPageBreak marker.
Replace-ParaXml "*Improve export forwarding code*" (
  '<w:p ' + $w + '>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Improve export forwarding code to detect and handle forward-by-</w:t></w:r>' +
    '<w:r><w:t>ordinal</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> explicitly rather than forcing the user to detect it and do string manipulation and conversion. </w:t></w:r>' +
  '</w:p>'
)
